# Add season record columns (Wins, Losses, Ties) to the DET_2008 sheet.
# The sheet currently spans A1:AC48 (header row 1, data rows 2-48).
# We append three new columns: AD = Wins, AE = Losses, AF = Ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of row 1, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) shares the same season record for this team.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 74  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 88  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
